# Applies the "3. osa tehtud" update to the "Nädal 4" sheet (time recording log).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 4")

# Row 7: Interruption Time (F7)
$ws.Range("F7").Value = 155

# Row 9: Interruption Time (F9)
$ws.Range("F9").Value = 110

# Row 10: Interruption Time (F10)
$ws.Range("F10").Value = 85

# Row 11: Interruption Time (F11) + clear comment (H11)
$ws.Range("F11").Value = 31
$ws.Range("H11").Value = ""

# Row 12: Stop time (D12), Interruption Time (F12), Activity (G12), clear Comments (H12)
$ws.Range("D12").Value = 0.72222222222222221
$ws.Range("F12").Value = 270
$ws.Range("G12").Value = "Kodutöö 3. osa + laadisin resharperi uuesti alla"
$ws.Range("H12").Value = ""

# Move the active selection to G15 (matches the author's cursor position after edits)
$ws.Range("G15").Select()

# Recalculate so the SUM(F7:F18) cached value picks up the new total
$excel.Calculate()
